# Case and Fatality Demographics Data Updated
# Weekly refresh (12.23.21 pull -> 12.30.21 pull): update the "Fatalities by
# Gender", "Fatalities by Age Group" and "Fatalities by Race-Ethnicity"
# worksheets with the new cumulative counts/percentages. The "Cases by *"
# worksheets are unaffected by this pull (their underlying numbers did not
# change), only the active-cell/selection bookkeeping moves as the author
# clicked through the tabs while reviewing the refreshed report.

$wb = $excel.ActiveWorkbook

$wsCasesAge     = $wb.Worksheets.Item("Cases by Age Group")
$wsCasesGender  = $wb.Worksheets.Item("Cases by Gender")
$wsCasesRace    = $wb.Worksheets.Item("Cases by RaceEthnicity")
$wsFatGender    = $wb.Worksheets.Item("Fatalities by Gender")
$wsFatAge       = $wb.Worksheets.Item("Fatalities by Age Group")
$wsFatRace      = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")

# ---------------------------------------------------------------------
# Fatalities by Gender: B/C were plain (non-formula) values, so both the
# counts and the already-computed percentages need to be re-entered.
# ---------------------------------------------------------------------

# Pull formatting for the refreshed cells from a sibling table that already
# uses the "plain data row" / "grand total row" styles, matching what
# happens when the new pull is pasted in over the old numbers.
$wsCasesAge.Range("B2:C4").Copy()
$wsFatGender.Range("B2:C4").PasteSpecial(-4122)
$wsCasesAge.Range("B15").Copy()
$wsFatGender.Range("B5").PasteSpecial(-4122)
$wsCasesAge.Range("C15").Copy()
$wsFatGender.Range("C5").PasteSpecial(-4122)

$wsFatGender.Range("B2").Value = 31233
$wsFatGender.Range("C2").Value = 31233/74578
$wsFatGender.Range("B3").Value = 43344
$wsFatGender.Range("C3").Value = 43344/74578
$wsFatGender.Range("B4").Value = 1
$wsFatGender.Range("C4").Value = 1/74578
$wsFatGender.Range("B5").Value = 74578
$wsFatGender.Range("C5").Value = 1

# The grand-total row no longer carries the old custom row height.
$wsFatGender.Rows.Item(5).AutoFit()

# ---------------------------------------------------------------------
# Fatalities by Age Group: B is plain values (C is a shared formula that
# recalculates automatically once B / B15 change).
# ---------------------------------------------------------------------

$wsCasesAge.Range("B2").Copy()
$wsFatAge.Range("B2:B14").PasteSpecial(-4122)
$wsCasesAge.Range("B15").Copy()
$wsFatAge.Range("B15").PasteSpecial(-4122)
$wsCasesAge.Range("C15").Copy()
$wsFatAge.Range("C15").PasteSpecial(-4122)

$wsFatAge.Range("B3").Value = 24
$wsFatAge.Range("B5").Value = 690
$wsFatAge.Range("B6").Value = 2206
$wsFatAge.Range("B7").Value = 5326
$wsFatAge.Range("B8").Value = 10098
$wsFatAge.Range("B9").Value = 7658
$wsFatAge.Range("B10").Value = 8976
$wsFatAge.Range("B11").Value = 9498
$wsFatAge.Range("B12").Value = 8969
$wsFatAge.Range("B13").Value = 21037
$wsFatAge.Range("B15").Value = 74578

# The grand-total row no longer carries the old custom row height.
$wsFatAge.Rows.Item(15).AutoFit()

# ---------------------------------------------------------------------
# Fatalities by Race-Ethnicity: B8/C8 are SUM/ratio formulas and pick up
# the new grand total automatically once B2:B7 are refreshed.
# ---------------------------------------------------------------------

$wsFatRace.Range("B2").Value = 1362
$wsFatRace.Range("B3").Value = 7897
$wsFatRace.Range("B4").Value = 32459
$wsFatRace.Range("B5").Value = 449
$wsFatRace.Range("B6").Value = 32366
$wsFatRace.Range("B7").Value = 45

# ---------------------------------------------------------------------
# Selection/active-tab bookkeeping as the author tabbed through the
# workbook reviewing the refresh (Cases sheets unaffected otherwise).
# ---------------------------------------------------------------------

$wsCasesAge.Activate()
$wsCasesAge.Range("B2:B14").Select()

$wsCasesGender.Range("B2:B4").Select()

$wsCasesRace.Range("B2:B7").Select()

$wsFatGender.Range("B5").Select()

$wsFatAge.Range("C26").Select()

$wsFatRace.Range("C7").Select()

$wsCasesAge.Activate()
